$d = $word.ActiveDocument

$replacements = @(
    @("51×86=", "71×60="),
    @("71×41=", "79×68="),
    @("31×66=", "78×76="),
    @("86×17=", "19×39="),
    @("43×23=", "56×78="),
    @("82×49=", "34×34="),
    @("50×13=", "47×11="),
    @("93×12=", "28×29="),
    @("83×64=", "51×53="),
    @("91×94=", "27×64="),
    @("45×38=", "46×58="),
    @("24×91=", "39×82="),
    @("41×59=", "54×32="),
    @("53×65=", "18×15="),
    @("66×20=", "67×57="),
    @("87×55=", "86×94="),
    @("15×34=", "13×54="),
    @("69×16=", "78×83="),
    @("39×19=", "43×43="),
    @("88×89=", "76×84="),
    @("81×58=", "33×29="),
    @("31×35=", "33×34="),
    @("84×48=", "95×24="),
    @("52×34=", "35×59="),
    @("73×62=", "99×64=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
